$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update APELLIDOS_Y_NOMBRES (column B) values, in the exact order
# that reproduces the target shared-string table ordering: rows 3-14
# (first occurrence of each new name) first, then row 2 (whose name is
# only reused there), then the remaining repeat rows. ---
$ws.Range("B3").Value = "Albert Dayhan Diaz"
$ws.Range("B4").Value = "Duvan Gutierrez Lobo"
$ws.Range("B5").Value = "Daniela Guzman Perez"
$ws.Range("B6").Value = "Perez Carmen Andrea"
$ws.Range("B7").Value = "Angulo Juan Camilo"
$ws.Range("B8").Value = "Cristian Olivar Isaza"
$ws.Range("B9").Value = "Thania Milena Perez"
$ws.Range("B10").Value = "Marlene Ballena Guzman"
$ws.Range("B11").Value = "Jose Sierra Guzman"
$ws.Range("B12").Value = "Diego Silva Benavides"
$ws.Range("B13").Value = "Camilo Andres Daza"
$ws.Range("B14").Value = "Jose Castellano Endry"
$ws.Range("B2").Value = "Elvis Galvis Galvis"

$ws.Range("B15").Value = "Albert Dayhan Diaz"
$ws.Range("B16").Value = "Duvan Gutierrez Lobo"
$ws.Range("B17").Value = "Daniela Guzman Perez"
$ws.Range("B18").Value = "Perez Carmen Andrea"
$ws.Range("B19").Value = "Angulo Juan Camilo"
$ws.Range("B20").Value = "Cristian Olivar Isaza"
$ws.Range("B21").Value = "Thania Milena Perez"
$ws.Range("B22").Value = "Marlene Ballena Guzman"
$ws.Range("B23").Value = "Jose Sierra Guzman"
$ws.Range("B24").Value = "Diego Silva Benavides"
$ws.Range("B25").Value = "Camilo Andres Daza"
$ws.Range("B26").Value = "Jose Castellano Endry"
$ws.Range("B27").Value = "Albert Dayhan Diaz"
$ws.Range("B28").Value = "Duvan Gutierrez Lobo"
$ws.Range("B29").Value = "Daniela Guzman Perez"
$ws.Range("B30").Value = "Perez Carmen Andrea"

# --- Numeric grade / fault-count corrections ---
$ws.Range("E2").Value = 3.5
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("F3").Value = 8
$ws.Range("F5").Value = 7
$ws.Range("H5").Value = 8

# --- Leftover "ghost" style (applyFont/applyAlignment, default look) on
#     the repeated-name cells, matching cellXfs index 4 in the target. ---
$ws.Range("B3").WrapText = $false
$ws.Range("B4").WrapText = $false
$ws.Range("B15").WrapText = $false
$ws.Range("B16").WrapText = $false
$ws.Range("B27").WrapText = $false
$ws.Range("B28").WrapText = $false

# --- Column B width ---
$ws.Columns.Item(2).ColumnWidth = 47.5

# --- Selection / scroll position ---
$ws.Range("C10").Select()
